$d = $word.ActiveDocument

# --- Question 1: add _GoBack bookmark, merge br+text runs ---
$rngQ1 = $d.Content
$rngQ1.Find.Execute("Question 1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rngQ1.Find.Found) { throw "Not found: Question 1" }
$paraQ1 = $rngQ1.Paragraphs(1)
$targetQ1 = $d.Range($paraQ1.Range.Start, $paraQ1.Range.End - 1)
$xmlQ1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>Question 1</w:t></w:r><w:r><w:br/><w:t>An individual’s tax price is defined as</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$targetQ1.InsertXML($xmlQ1)

# --- Question 2: remove _GoBack bookmark, merge br+text runs ---
$rngQ2 = $d.Content
$rngQ2.Find.Execute("Question 2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rngQ2.Find.Found) { throw "Not found: Question 2" }
$paraQ2 = $rngQ2.Paragraphs(1)
$targetQ2 = $d.Range($paraQ2.Range.Start, $paraQ2.Range.End - 1)
$xmlQ2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t>Question 2</w:t></w:r><w:r><w:br/><w:t>The median voter is powerful because</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$targetQ2.InsertXML($xmlQ2)

# --- Question 3: merge br+text runs ---
$rngQ3 = $d.Content
$rngQ3.Find.Execute("Question 3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rngQ3.Find.Found) { throw "Not found: Question 3" }
$paraQ3 = $rngQ3.Paragraphs(1)
$targetQ3 = $d.Range($paraQ3.Range.Start, $paraQ3.Range.End - 1)
$xmlQ3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t>Question 3</w:t></w:r><w:r><w:br/><w:t>Inefficient public operations can persist because</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$targetQ3.InsertXML($xmlQ3)

# --- Question 4: merge br+text runs ---
$rngQ4 = $d.Content
$rngQ4.Find.Execute("Question 4", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rngQ4.Find.Found) { throw "Not found: Question 4" }
$paraQ4 = $rngQ4.Paragraphs(1)
$targetQ4 = $d.Range($paraQ4.Range.Start, $paraQ4.Range.End - 1)
$xmlQ4 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t>Question 4</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">A </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Lindahl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> equilibrium</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$targetQ4.InsertXML($xmlQ4)

# --- Question 5: merge br+text runs ---
$rngQ5 = $d.Content
$rngQ5.Find.Execute("Question 5", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rngQ5.Find.Found) { throw "Not found: Question 5" }
$paraQ5 = $rngQ5.Paragraphs(1)
$targetQ5 = $d.Range($paraQ5.Range.Start, $paraQ5.Range.End - 1)
$xmlQ5 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t>Question 5</w:t></w:r><w:r><w:br/><w:t>Arrow’s impossibility theorem shows that no government decision rule can always</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$targetQ5.InsertXML($xmlQ5)

# --- "sets supply equal to demand." paragraph: color red + append page cite ---
$rngColor = $d.Content
$rngColor.Find.Execute("sets supply equal to demand.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rngColor.Find.Found) { throw "Not found: sets supply equal to demand." }
$paraColor = $rngColor.Paragraphs(1)
$paraColor.Range.Font.Color = 255

$rngSS = $d.Content
$rngSS.Find.Execute("sets supply equal to demand.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rngSS.Find.Found) { throw "Not found: sets supply equal to demand." }
$paraSS = $rngSS.Paragraphs(1)
$targetSS = $d.Range($paraSS.Range.Start, $paraSS.Range.End - 1)
$xmlSS = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>sets supply equal to demand.</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> (p. 249, 2nd paragraph of section)</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$targetSS.InsertXML($xmlSS)

# --- Insert Answer Explanation block after "assigns tax prices..." paragraph ---
$rngAns = $d.Content
$rngAns.Find.Execute("assigns tax prices proportional to average benefits.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rngAns.Find.Found) { throw "Not found: assigns tax prices..." }
$paraAns = $rngAns.Paragraphs(1)
$insertPoint = $d.Range($paraAns.Range.End - 1, $paraAns.Range.End - 1)
$xmlAns = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>Answer Explanation:</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Lindahl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> equilibrium equates sum of marginal benefits (</w:t></w:r><w:r><w:t xml:space="preserve">i.e., aggregate marginal benefit, </w:t></w:r><w:r><w:t>not individual marginal benefit) to the marginal cost of production (p. 251).</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Lindahl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> equilibrium is a set of tax prices that add up to the marginal cost of production</w:t></w:r><w:r><w:t>.  Aggregate marginal benefits</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>denote the total marginal willingness to pay</w:t></w:r><w:r><w:t>, which is the same as the sum of tax prices</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>(p. 251).</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$insertPoint.InsertXML($xmlAns)

Write-Host "All edits applied."